$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 337, shifting existing rows 337:416 down to 338:417
$ws.Rows.Item(337).Insert()

# Populate the new row 337 with data (copy constant columns from the row below, then set changed values)
$ws.Range("A337").Value = 3
$ws.Range("B337").Value = "Femacal de La Calera"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value = 44782
$ws.Range("D337").NumberFormat = $ws.Range("D338").NumberFormat
$ws.Range("E337").Value = 5
$ws.Range("F337").Value = 100112031
$ws.Range("G337").Value = "Poroto verde"
$ws.Range("H337").Value = "Magnum"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 85
$ws.Range("K337").Value = 33000
$ws.Range("L337").Value = 34000
$ws.Range("M337").Value = 33529
$ws.Range("N337").Value = "$/malla 25 kilos"
$ws.Range("O337").Value = "Región de Arica y Parinacota"
$ws.Range("P337").Value = 1341
$ws.Range("Q337").Value = 25
$ws.Range("R337").Value = "Hortaliza"
